# Add new "shallow" scenario plus temp/zoop variants for the Swim-Advect-Diffuse
# block, and two new result blocks ("Swim-Advect" and "Swim") on the summary
# sheet, mirroring the existing "deep" row's layout/format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary")

# ---- Helper data: the six scenario labels (in on-sheet order) and the three
# blocks of F-column results that go with them. ----
$labels = @(
    "velocity, daily 2000, bio=100, dt=1 hr, deep",
    "velocity, daily 2000, bio=100, dt=1 hr, shallow",
    "velocity, daily 2000, bio=100, dt=1 hr, temp const",
    "velocity, daily 2000, bio=100, dt=1 hr, temp daily",
    "velocity, daily 2000, bio=100, dt=1 hr, zoop const",
    "velocity, daily 2000, bio=100, dt=1 hr, zoop daily"
)

# Swim-Advect-Diffuse block (rows 65-70) -- header already exists in A64
$valsSwimAdvectDiffuse = @(0.1191, -0.1785, 12.386, 0.8182, 0.2999, 0.3306)

# Swim-Advect block (rows 72-78)
$valsSwimAdvect = @(0.6807, -0.0136, 90.7487, 3.2086, 0.9887, 0.6388)

# Swim block (rows 80-86)
$valsSwim = @(75.871, 0.0142, 276.3927, 3.2788, 276.5259, 132.6202)

function Write-ScenarioBlock($startRow, $values) {
    for ($i = 0; $i -lt $labels.Count; $i++) {
        $r = $startRow + $i
        $ws.Cells.Item($r, 1).Value = $labels[$i]
        $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 5)).NumberFormat = "0.00E+00"
        $f = $ws.Cells.Item($r, 6)
        $f.NumberFormat = "0.00E+00"
        $f.Value = $values[$i]
    }
}

# --- Block 1: continues the existing "Swim-Advect-Diffuse" section ---
Write-ScenarioBlock 65 $valsSwimAdvectDiffuse

# Blank spacer row 71 (style-only, matches the blank rows elsewhere)
$ws.Range("B71:F71").NumberFormat = "0.00E+00"

# --- Block 2: new "Swim-Advect" section header (row 72) + data (73-78) ---
$ws.Cells.Item(72, 1).Value = "Swim-Advect"
$ws.Cells.Item(72, 1).Font.Bold = $true
$ws.Range("B72:F72").NumberFormat = "0.00E+00"

Write-ScenarioBlock 73 $valsSwimAdvect

# Blank spacer row 79
$ws.Range("B79:F79").NumberFormat = "0.00E+00"

# --- Block 3: new "Swim" section header (row 80) + data (81-86) ---
$ws.Cells.Item(80, 1).Value = "Swim"
$ws.Cells.Item(80, 1).Font.Bold = $true
$ws.Range("B80:F80").NumberFormat = "0.00E+00"

Write-ScenarioBlock 81 $valsSwim

# ---- Best-effort view state (scroll position / selection) ----
$ws.Activate()
$ws.Range("F84").Select()
$excel.ActiveWindow.ScrollRow = 52
